# Applies corrected Diebold-Mariano test statistics (DM_Stat, column C)
# and p-values (P_Value, column D) for rows 2-11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.7910369926831037
$ws.Range("D2").Value = 0.4373678028254557

$ws.Range("C3").Value = -0.1910957353320806
$ws.Range("D3").Value = 0.8502031780747414

$ws.Range("C4").Value = 0.04522190003856343
$ws.Range("D4").Value = 0.9643384124858849

$ws.Range("C5").Value = -0.6633795859031241
$ws.Range("D5").Value = 0.5139793536714663

$ws.Range("C6").Value = -0.7827676137257374
$ws.Range("D6").Value = 0.4421085662956872

$ws.Range("C7").Value = -0.297281874158327
$ws.Range("D7").Value = 0.7690402375463683

$ws.Range("C8").Value = -1.082166271550063
$ws.Range("D8").Value = 0.2908980649845574

$ws.Range("C9").Value = 0.1562304203015938
$ws.Range("D9").Value = 0.8772761606385766

$ws.Range("C10").Value = -0.4484750573824606
$ws.Range("D10").Value = 0.65819632758061

$ws.Range("C11").Value = -0.4077096691111827
$ws.Range("D11").Value = 0.687427203291465
